$d = $word.ActiveDocument

# --- Change 1: merge "A " + "assistente social " + "fez a rotina de
# autenticação no sistema conforme " into a single run with the same
# formatting (sz 28). A plain Find/Replace over the exact combined text
# collapses the three adjacent same-format runs into one. ---
$d.Content.Find.Execute(
    "A assistente social fez a rotina de autenticação no sistema conforme ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "A assistente social fez a rotina de autenticação no sistema conforme ",
    2) | Out-Null

# --- Change 2: "Tela 02" -> "Tela " + "0100" (the digits becoming their
# own run, keeping the same red/size formatting as before). Turning on
# track-changes for the edit keeps the insertion as a distinct run
# instead of Word silently re-merging identical-format neighbours;
# accepting the revision afterwards leaves the two runs behind exactly
# as authored. ---
$d.TrackRevisions = $true

$r = $d.Content
$r.Find.Execute("Tela 02", $false, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$numRange = $d.Range($r.Start + 5, $r.End)
$numRange.Text = "0100"

$d.TrackRevisions = $false

# Accept only the revisions we just made (rather than AcceptAllRevisions,
# which also touches unrelated parts of the document).
for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
    $d.Revisions.Item($i).Accept()
}
